# "Generate Report for Handback" - refresh the localization-status report
# after a successful handback: the status moves from "Ready for handoff" to
# "Handed back: in sync with en-US" everywhere it appears, the per-language
# "Latest Handback DateTime" is stamped with the new handback time, the
# (now resolved) "Error Detail" about a stale handback file is cleared, and
# the columns that display that information are widened/narrowed to fit.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: zh-cn / de-de status columns -------------------------
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus

# --- zh-cn sheet ------------------------------------------------------------
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("K2").Value = "2016-08-15 10:47:55"
$zhcn.Range("P2").Value = ""

# --- de-de sheet ------------------------------------------------------------
$dede.Range("C2").Value = $newStatus
$dede.Range("K2").Value = "2016-08-15 10:48:04"
$dede.Range("P2").Value = ""

# --- Column width refresh (status/error columns got wider/narrower text) ---
$overview.Columns.Item(5).ColumnWidth = 29.166666666666668
$overview.Columns.Item(6).ColumnWidth = 29.166666666666668

$zhcn.Columns.Item(3).ColumnWidth  = 29.166666666666668
$zhcn.Columns.Item(16).ColumnWidth = 12.833333333333334

$dede.Columns.Item(3).ColumnWidth  = 29.166666666666668
$dede.Columns.Item(16).ColumnWidth = 12.833333333333334
